$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new "missing data" rows (74-76) under the existing
# running-hours table, continuing the running total formula in column C. ---

# New dates (column A), carrying over the date number-format/style used
# by the row directly above (row 73) so no new style entry is created.
$ws.Range("A74").Value = 45629
$ws.Range("A75").Value = 45643
$ws.Range("A76").Value = 45644

$ws.Range("A73").Copy() | Out-Null
$ws.Range("A74:A76").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Hours worked that day (column B).
$ws.Range("B74").Value = 8
$ws.Range("B75").Value = 7
$ws.Range("B76").Value = 3

# Running total (column C) - continues the C70+B71 style accumulation.
$ws.Range("C74").Formula = "=C73+B74"
$ws.Range("C75").Formula = "=C74+B75"
$ws.Range("C76").Formula = "=C75+B76"

# Match the workbook author's final selection/active cell (C76, the last
# filled cell in the running-total column).
$ws.Range("C76").Select() | Out-Null
